$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Fill in the new "XXXX" placeholder answer in column U ("Can we contact
#    you ...") for every existing response row that previously had this cell
#    blank. The cell already carries style s="1" (text style), so setting
#    .Value reuses that style and adds "XXXX" to the shared-strings table.
# ---------------------------------------------------------------------------
$rowsToFill = @(4, 5, 6, 7, 10, 13, 15, 18, 22, 23, 24, 25, 26)
foreach ($r in $rowsToFill) {
    $ws.Cells.Item($r, 21).Value = "XXXX"
}

# ---------------------------------------------------------------------------
# 2) Append the new survey response as row 27.
#    Formats are copied cell-by-cell from the row above (row 26) so the
#    existing shared cellXf indices (s="1" for text, s="2" for the Timestamp
#    date) are reused instead of new styles being created; values are then
#    written on top of the copied formatting.
# ---------------------------------------------------------------------------
$newRow = 27
$templateRow = 26

$filledCols = @(1, 2, 3, 4, 5, 6, 7, 8, 10, 12, 14, 16, 18, 20, 21)
foreach ($c in $filledCols) {
    $ws.Cells.Item($templateRow, $c).Copy()
    $ws.Cells.Item($newRow, $c).PasteSpecial(-4122)
}

$ws.Cells.Item($newRow, 1).Value = 45141.428112534719   # Timestamp
$ws.Cells.Item($newRow, 2).Value = "Europe"
$ws.Cells.Item($newRow, 3).Value = "Enterprise/Company"
$ws.Cells.Item($newRow, 4).Value = "Developer"
$ws.Cells.Item($newRow, 5).Value = "Medium (51- 250 employees)"
$ws.Cells.Item($newRow, 6).Value = "2 - 5 years"
$ws.Cells.Item($newRow, 7).Value = "2 - 5 projects"
$ws.Cells.Item($newRow, 8).Value = "Yes"
$ws.Cells.Item($newRow, 10).Value = 4
$ws.Cells.Item($newRow, 12).Value = 4
$ws.Cells.Item($newRow, 14).Value = 4
$ws.Cells.Item($newRow, 16).Value = 2
$ws.Cells.Item($newRow, 18).Value = 4
$ws.Cells.Item($newRow, 20).Value = "Yes"
$ws.Cells.Item($newRow, 21).Value = "XXXX"

# Keep the new row's height consistent with every other data row.
$ws.Rows.Item($newRow).RowHeight = 15.75

# ---------------------------------------------------------------------------
# 3) Move the active selection, matching the author's last-saved cursor
#    position after adding the row above.
# ---------------------------------------------------------------------------
[void]$ws.Range("U31").Select()
